$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column C header from 'Answer' to 'Context', and add new column D 'Answer'
$ws.Range("C1").Value = 'Context'
$ws.Range("D1").Value = 'Answer'

# New short answers for column D, rows 2-61 (in row order)
$answers = @(
  'There might be no significant differences in conscientious intelligence sub-dimensions between first and fourth-year nursing students because the study found that factors such as willingness to choose nursing as a profession, gender (female), and not perceiving nursing as solely a means to earn money were more influential in determining conscientious intelligence scores than academic year level.',
  'Yes, according to the study, nursing students who do not perceive nursing as solely a means to earn money demonstrate higher levels of conscientious intelligence.',
  'b) Willingness to choose the profession',
  'The study quantitatively demonstrates the relationship between IQ and brain resilience by utilizing resting-state fMRI graph-theoretical analysis in 102 healthy individuals. They observe enhanced brain robustness to targeted attacks (TA) in individuals with higher IQ, indicating increased distributed processing capacity despite the systematic loss of important nodes. Additionally, they identify specific neocortical regions, primarily involved in language and memory processing networks, associated with brain resilience in higher IQ individuals, while emotional processing regions are more relevant for lower IQ individuals.',
  'Yes, individuals with higher IQ show resilience to both targeted and random attacks on the brain, as evidenced by the study''s findings. The enhanced brain robustness to targeted attacks (TA) in individuals with higher IQ indicates increased distributed processing capacity despite the systematic loss of important nodes. Additionally, brain resilience in higher IQ individuals is supported by a set of neocortical regions mainly involved in language and memory processing networks.',
  'a) Higher IQ: Neocortical regions mainly belonging to language and memory processing networks.
b) Lower IQ: Regions related to emotional processing.',
  'Exposure to lead-based industries can negatively impact the intelligence of school-aged children by interfering with brain development and function, leading to cognitive impairments, as evidenced by a negative correlation between blood lead levels and intelligence scores.',
  'Nutritional supplementation and hygiene education had a greater impact on reducing lead levels in the experimental group compared to improving intelligence scores.',
  'c) Lead exposure may not directly impact hemoglobin levels.',
  'The study provides evidence supporting the use of new tools for measuring mindsets in medical education through the satisfactory psychometric properties demonstrated by the new items. Additionally, the study''s findings reveal the impact of mindsets on medical education, suggesting the potential usefulness of these tools in understanding and assessing clinical supervisors'' beliefs about the fixed versus learnable nature of attributes such as intelligence, clinical reasoning, moral character, and empathy.',
  'The greatest percentage of clinical supervisors viewed empathy as fixed.',
  'A clinical supervisor''s mindset about a particular attribute might impact their teaching practices when they believe that attribute to be fixed rather than learnable. For example, if a supervisor believes that empathy is fixed and cannot be developed through education and training, they may be less likely to emphasize or prioritize empathy-related skills in their teaching. This could result in a teaching approach that focuses less on cultivating empathy in medical students or residents, potentially affecting the quality of patient care and interpersonal interactions.',
  'Medical educators, curriculum developers, and policymakers might benefit from better understanding the impact of mindsets on medical education.',
  'The researchers measured the clinical supervisors'' mindsets using surveys containing existing instruments for intelligence and moral character, as well as 18 new items for assessing clinical reasoning and empathy.',
  'The PCTB assesses areas of cognition including physical cognition (e.g., quantity discrimination) and social cognition (e.g., gaze following) in chimpanzees.',
  'The factor in the PCTB most strongly related to HDT efficiency was chimpanzee g scores.',
  'A chimpanzee''s self-control abilities might be most critical for successful performance on the PCTB when they need to sustain delay of gratification while engaging in tasks that require cognitive processing, such as quantity discrimination or gaze following.',
  'Besides chimpanzees, humans also demonstrate a link between general intelligence and self-control behavior.',
  'The Hybrid Delay Task (HDT) measures self-control in chimpanzees by assessing their inter-temporal choices and capacity for sustained delay of gratification.',
  'The NART-derived WAIS-IV FSIQ estimates are approximately equidistant from the highly discrepant WAIS (low) and WAIS-R (high) values, particularly at the lower end of the distribution.',
  'The pattern of correlations between NART performance and the different WAIS-IV measures showed strong correlations with WAIS-IV Full-Scale IQ (FSIQ) and more moderate correlations with constituent index scores. There were differences in the strength of correlations across the board, with FSIQ showing stronger correlations compared to constituent index scores. Additionally, correlations with NART error scores were more moderate compared to correlations with FSIQ.',
  'The most anterior face patch (AM) might show a weaker representation of non-face stimuli compared to other face patches because it likely specializes more in processing facial information rather than non-face visual stimuli.',
  'Researchers use population decoding to characterize the neural activity within face patches by analyzing patterns of activity across neurons to decode information related to discriminating between face and nonface objects, individual faces, and nonface stimuli.',
  'Yes, according to the study, patients with schizophrenia who have low IQ and severe negative symptoms are more likely to demonstrate insufficient effort on neuropsychological testing, even if positive symptoms are less pronounced.',
  'Yes, according to the study, an increased prevalence of "deficit schizophrenia" within a clinical population might lead to a corresponding increase in insufficient effort during neuropsychological evaluations.',
  'Besides personality traits and intelligence, other factors that might be used as dimensions within an individual''s "Mind-space" could include beliefs, emotions, values, attitudes, preferences, and social roles.',
  'Yes, an individual''s "Mind-space" model can change or evolve over time as they gain new experiences, encounter different social situations, and acquire new knowledge about the variability of mental states and personality traits.',
  'b) The accuracy of an individual''s model of how minds covary',
  'The revised Carter model, incorporating catastrophes, suggests that past estimates of critical steps in the evolution of intelligent life may have been underestimated. It proposes that catastrophes can accelerate or cluster critical steps, potentially leading to a more rapid emergence of intelligent life.',
  'Yes, the apparent increase in biological complexity over the past 500 million years could be partially explained by previously overlooked evolutionary transitions, including those influenced by catastrophes, according to the generalized Carter model.',
  'The neo-catastrophes hypothesis provides a more plausible explanation for the Fermi paradox, which observes the rarity of intelligent life in the observable Universe.',
  'Besides drug-use impulsivity and intelligence, other factors that might contribute to differences in heroin delay discounting (DD) during states of satiation and withdrawal could include individual differences in physiological arousal, emotional state, or cognitive functioning specific to each state. Additionally, environmental factors, such as the presence of cues associated with drug use or stressors, could also influence DD during these states.',
  'No, the factors predicting heroin delay discounting (DD) in a state of satiation partially overlap with those predicting DD in a state of withdrawal. Lower intelligence and higher IRQ Speed (to return to drug use) predicted greater discounting in both states, while lower IRQ Capacity for Delay predicted greater discounting only in the satiation condition.',
  'c) Heroin delay discounting is influenced by both pharmacological state and individual characteristics.',
  'The smart home system utilizes communication technology, Internet connection technology, automatic fire control technology, network wiring technology, and visual and audio transmission technology.',
  'The web server intelligent management system offers remote access for controlling home devices.',
  'A user might choose to control their smart home devices from a mobile terminal when they are on the go or prefer a more portable option. They might opt for a web interface when they are at a stationary location with access to a computer or prefer a larger screen for control.',
  'The intended users of this smart home system are homeowners or individuals looking to manage and control household appliances and systems for convenience and efficiency.',
  'The system enables control of lighting, air conditioning, and other appliances wirelessly through mobile terminals or web interfaces connected to the smart home management system.',
  'The purpose of integrating space topological cognition into the model is to serve as a path planning module, enabling the model to obtain a more generalized approach to navigation.',
  'The space topological cognition component is responsible for guiding detour behavior in the navigation architecture.',
  'An agent would rely on space topological cognition instead of direct locomotion for navigation when a path to the target is unavailable, requiring detour behavior.',
  'The curiosity-driven rewards influence the agent''s exploration behavior by incentivizing it to explore novel and uncertain regions of the environment, leading to more comprehensive exploration and learning.',
  'The emerging field of study addressing the limitations of traditional safety ergonomics in the era of big data is called Intelligent Safety Ergonomics (ISE).',
  'The integration of information technology and safety science, as seen in Intelligent Safety Ergonomics (ISE), can lead to a more holistic understanding of human-machine-environment systems by providing deeper insights into their interactions and dynamics. This understanding can lead to improved workplace safety and sustainability through more informed decision-making and the development of innovative solutions tailored to specific contexts and challenges.',
  'The study challenges previous findings by demonstrating that awareness of the location-reward association, rather than goal-driven manipulations, is crucial for inducing reward biases to high-reward locations.',
  'Yes, the researchers found a significant location-reward effect even when participants were not actively searching for high-reward locations.',
  'b) Awareness of location-reward association',
  'The primary goals of the Hungarian Primary Care Model Programme were to increase access to primary care for disadvantaged population groups and to offer new preventive services for all clients.',
  'The primary beneficiaries of the health mediator''s work were the local communities, including vulnerable Roma members, who aimed to improve access to health services.',
  'The health mediators'' work increased their self-esteem and social status by allowing them to become equal members of the primary care team, contributing to the improvement of access to health services for disadvantaged populations.',
  'The use of health mediators is particularly important in areas with disadvantaged Roma populations because they can improve access to health services for this vulnerable group, addressing their specific needs and increasing their utilization of preventive and primary care.',
  'The GP clusters implementing the Primary Care Model Programme were located in Hungary.',
  'The evaluation of the Primary Care Model Programme took place between 2013 and 2017 in Hungary.',
  'The health mediators'' job aspects related to their personal characteristics, such as motivation, increased self-esteem, social status, and health consciousness, received the most positive feedback.',
  'Clinicians can improve their communication by understanding patients'' values, treating them with respect, acknowledging emotional and social realities, recognizing implicit biases, being mindful of language around weight and food, and using verbal and nonverbal approaches to convey empathic concern.',
  'It''s crucial for clinicians to recognize their own implicit biases about weight and food to ensure they communicate with patients in a respectful and nonjudgmental manner, fostering stronger relationships and promoting better adherence to dietary recommendations.',
  'A leader''s level of trait empathy influences their emotional reactions and perceived effectiveness after providing negative feedback. Leaders higher in trait empathy report feeling less attentive and more distressed, while those lower in trait empathy feel more attentive and less distressed.',
  'Yes, the study suggests that there are potential downsides to leader empathy in the context of giving negative feedback.',
  'Nurse leaders can use emotional intelligence to specifically address the challenges of stress, exhaustion, and the risk of moral injury brought on by the COVID-19 pandemic by supporting their teams, fostering effective personal relationships, and directing them towards common goals while demonstrating sensitivity to their psychological health and well-being.'
)

for ($i = 0; $i -lt $answers.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 4).Value = $answers[$i]
}
